$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 46.29121633333333
$ws.Range("N2").Value = 138.873649
$ws.Range("O2").Value = 0.3133663986859022
$ws.Range("P2").Value = 0.3133663986859022
$ws.Range("Q2").Value = 15.55770628936111
$ws.Range("R2").Value = 140.01935660425
$ws.Range("S2").Value = 0.3133663986859022
$ws.Range("T2").Value = 0.3133663986859022

# Row 3
$ws.Range("M3").Value = 46.81622333333333
$ws.Range("O3").Value = 0.3169204109998198
$ws.Range("P3").Value = 0.3169204109998198
$ws.Range("S3").Value = 0.3169204109998198
$ws.Range("T3").Value = 0.3169204109998198

# Row 4
$ws.Range("M4").Value = 38.53544233333333
$ws.Range("N4").Value = 115.606327
$ws.Range("O4").Value = 0.2608640200510233
$ws.Range("P4").Value = 0.2608640200510233
$ws.Range("Q4").Value = 12.95111991086111
$ws.Range("R4").Value = 116.56007919775
$ws.Range("S4").Value = 0.2608640200510233
$ws.Range("T4").Value = 0.2608640200510233

# Row 5
$ws.Range("M5").Value = 16.07945366666667
$ws.Range("N5").Value = 48.238361
$ws.Range("O5").Value = 0.1088491702632547
$ws.Range("P5").Value = 0.1088491702632547
$ws.Range("Q5").Value = 5.404036386472221
$ws.Range("R5").Value = 48.63632747824999
$ws.Range("S5").Value = 0.1088491702632547
$ws.Range("T5").Value = 0.1088491702632547
